$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 12, 9),
    @(3, 12, 9),
    @(4, 12, 9),
    @(5, 12, 9),
    @(6, 12, 9),
    @(7, 11, 9),
    @(8, 13, 8),
    @(9, 12, 9),
    @(10, 12, 9),
    @(11, 12, 9),
    @(12, 12, 9),
    @(13, 12, 9),
    @(14, 12, 10),
    @(15, 12, 8),
    @(16, 11, 10)
)

foreach ($row in $values) {
    $r = $row[0]
    $b = $row[1]
    $c = $row[2]
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
}
